$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value.
# "Coin"/"Link" rows shifted down one rank (a new coin -- OKB -- entered the
# top 50 at #8), and every Price/Volume(1h) figure was refreshed from the feed.
$edits = @(
    ,@("D2", "28.541.34")
    ,@("E2", "  -1.32%  ")
    ,@("D3", "1.887.31")
    ,@("E3", "  +0.22%  ")
    ,@("E4", "  +0.36%  ")
    ,@("D5", "326.37")
    ,@("E5", "  -0.44%  ")
    ,@("E6", "  +0.40%  ")
    ,@("E7", "  -1.65%  ")
    ,@("E8", "  -2.49%  ")
    ,@("B9", "OKB")
    ,@("C9", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb")
    ,@("D9", "46.77")
    ,@("E9", "  +0.23%  ")
    ,@("B10", "Dogecoin")
    ,@("C10", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge")
    ,@("D10", "0.07867")
    ,@("E10", "  -0.91%  ")
    ,@("B11", "Polygon")
    ,@("C11", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic")
    ,@("D11", "0.9997")
    ,@("E11", "  +2.15%  ")
    ,@("B12", "Solana")
    ,@("C12", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol")
    ,@("D12", "21.61")
    ,@("E12", "  -3.81%  ")
    ,@("B13", "WrappedEther")
    ,@("C13", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth")
    ,@("D13", "1.895.15")
    ,@("E13", "  +5.15%  ")
    ,@("B14", "Chainlink")
    ,@("C14", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link")
    ,@("D14", "7.068")
    ,@("E14", "  +0.61%  ")
    ,@("B15", "Polkadot")
    ,@("C15", "https://coinranking.com/coin/25W7FG7om+polkadot-dot")
    ,@("D15", "5.709")
    ,@("E15", "  -1.11%  ")
    ,@("B16", "TRON")
    ,@("C16", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx")
    ,@("D16", "0.06966")
    ,@("E16", "  -0.31%  ")
    ,@("B17", "Litecoin")
    ,@("C17", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc")
    ,@("D17", "87.51")
    ,@("E17", "  -1.69%  ")
    ,@("B18", "BinanceUSD")
    ,@("C18", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd")
    ,@("D18", "1.009")
    ,@("E18", "  +0.40%  ")
    ,@("B19", "ShibaInu")
    ,@("C19", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib")
    ,@("D19", "0.00001003")
    ,@("E19", "  -1.23%  ")
    ,@("B20", "Avalanche")
    ,@("C20", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax")
    ,@("D20", "17.17")
    ,@("E20", "  +0.60%  ")
    ,@("B21", "Dai")
    ,@("C21", "https://coinranking.com/coin/MoTuySvg7+dai-dai")
    ,@("D21", "1.006")
    ,@("E21", "  +0.25%  ")
    ,@("B22", "WrappedBTC")
    ,@("C22", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc")
    ,@("D22", "28.579.99")
    ,@("E22", "  -1.16%  ")
    ,@("B23", "Uniswap")
    ,@("C23", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni")
    ,@("D23", "5.321")
    ,@("E23", "  -0.95%  ")
    ,@("B24", "Cosmos")
    ,@("C24", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom")
    ,@("D24", "10.98")
    ,@("E24", "  -1.50%  ")
    ,@("B25", "WrappedliquidstakedEther2.0")
    ,@("C25", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth")
    ,@("D25", "2.134.56")
    ,@("E25", "  +4.62%  ")
    ,@("B26", "Toncoin")
    ,@("C26", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton")
    ,@("D26", "2.063")
    ,@("E26", "  -2.70%  ")
    ,@("B27", "Monero")
    ,@("C27", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr")
    ,@("D27", "154.86")
    ,@("E27", "  +0.78%  ")
    ,@("B28", "EthereumClassic")
    ,@("C28", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc")
    ,@("D28", "19.37")
    ,@("E28", "  -0.59%  ")
    ,@("B29", "InternetComputer(DFINITY)")
    ,@("C29", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp")
    ,@("D29", "5.817")
    ,@("E29", "  +0.38%  ")
    ,@("B30", "LidoDAOToken")
    ,@("C30", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo")
    ,@("D30", "1.956")
    ,@("E30", "  -2.97%  ")
    ,@("B31", "BitcoinCash")
    ,@("C31", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch")
    ,@("D31", "118.22")
    ,@("E31", "  -1.57%  ")
    ,@("B32", "Stellar")
    ,@("C32", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm")
    ,@("D32", "0.09322")
    ,@("E32", "  -0.90%  ")
    ,@("B33", "ImmutableX")
    ,@("C33", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx")
    ,@("D33", "0.9226")
    ,@("E33", "  -2.53%  ")
    ,@("B34", "Filecoin")
    ,@("C34", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil")
    ,@("D34", "5.292")
    ,@("E34", "  -0.90%  ")
    ,@("B35", "ARBITRUM")
    ,@("C35", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb")
    ,@("D35", "1.333")
    ,@("E35", "  -1.70%  ")
    ,@("B36", "HuobiToken")
    ,@("C36", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht")
    ,@("D36", "3.271")
    ,@("E36", "  -2.44%  ")
    ,@("B37", "Hedera")
    ,@("C37", "https://coinranking.com/coin/jad286TjB+hedera-hbar")
    ,@("D37", "0.05760")
    ,@("E37", "  -3.02%  ")
    ,@("B38", "TrustWalletToken")
    ,@("C38", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt")
    ,@("D38", "1.156")
    ,@("E38", "  +0.42%  ")
    ,@("B39", "VeChain")
    ,@("C39", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet")
    ,@("D39", "0.02069")
    ,@("E39", "  -2.98%  ")
    ,@("B40", "FraxShare")
    ,@("C40", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs")
    ,@("D40", "7.795")
    ,@("E40", "  -2.16%  ")
    ,@("B41", "TheSandbox")
    ,@("C41", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand")
    ,@("D41", "0.5663")
    ,@("E41", "  -1.61%  ")
    ,@("B42", "Algorand")
    ,@("C42", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo")
    ,@("D42", "0.1788")
    ,@("E42", "  -0.71%  ")
    ,@("B43", "Aptos")
    ,@("C43", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt")
    ,@("D43", "9.747")
    ,@("E43", "  -2.86%  ")
    ,@("B44", "EnergySwap")
    ,@("C44", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens")
    ,@("D44", "11.75")
    ,@("E44", "  -1.56%  ")
    ,@("B45", "Cronos")
    ,@("C45", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro")
    ,@("D45", "0.07150")
    ,@("E45", "  -1.46%  ")
    ,@("B46", "Decentraland")
    ,@("C46", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana")
    ,@("D46", "0.5346")
    ,@("E46", "  -0.37%  ")
    ,@("B47", "RenderToken")
    ,@("C47", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr")
    ,@("D47", "2.181")
    ,@("E47", "  +1.34%  ")
    ,@("B48", "NEARProtocol")
    ,@("C48", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near")
    ,@("D48", "1.837")
    ,@("E48", "  -1.29%  ")
    ,@("B49", "WEMIXToken")
    ,@("C49", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix")
    ,@("D49", "1.116")
    ,@("E49", "  -2.06%  ")
    ,@("B50", "Quant")
    ,@("C50", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt")
    ,@("D50", "112.63")
    ,@("E50", "  -1.57%  ")
    ,@("B51", "MXToken")
    ,@("C51", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx")
    ,@("D51", "2.472")
    ,@("E51", "  +4.33%  ")
)

foreach ($edit in $edits) {
    $ref = $edit[0]
    $val = $edit[1]
    $cell = $ws.Range($ref)
    # Guard numeric-looking text (e.g. "326.37", "0.00001003") so Excel keeps
    # storing it as text instead of silently re-typing the cell as a number,
    # while leaving the cells original style/formatting untouched.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = $origStyle
}
